$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("Q1").Value = "Test"
